$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix B51: was stored as text "5", should be a real number 5 ---
$ws.Cells.Item(51, 2).Value = 5

# --- Append new row 52 (new annotation entry for Ying Tang) ---
$ws.Cells.Item(52, 1).Value = "Ying Tang"

# B52 must stay as a *text* "3" (matches the source data's inconsistency
# where this particular score was entered as a string). Assigning the
# string "3" directly would auto-coerce to a number, so we compute it as
# text via a formula and then paste the result back as a plain value -
# this keeps the cell's stored type as text without leaving a lingering
# number-format/style behind.
$ws.Cells.Item(52, 2).Formula = '=TEXT(3,"0")'
$ws.Cells.Item(52, 2).Copy()
$ws.Cells.Item(52, 2).PasteSpecial(-4163)

$ws.Cells.Item(52, 3).Value = "无"
$ws.Cells.Item(52, 4).Value = "DIS"
$ws.Cells.Item(52, 5).Value = "MET"
$ws.Cells.Item(52, 6).Value = "6dbc86e6-aac5-4bea-af0c-fc9177dfd16b"
$ws.Cells.Item(52, 7).Value = "BkJ3ibb0-_annotated.xlsx"
$ws.Cells.Item(52, 8).Value = "Furthermore, we have not optimized the running time of our algorithm, as it was not the focus of this work."
